{"js": "// Remove the paragraph that holds the italic \"Genesis\" run directly\n// following the \"GEN\" (Heading 2) paragraph. The surrounding, differently\n// styled \"Genesis\" Heading 2 paragraph (further below in the document) must\n// stay untouched, so we match on BOTH the exact text and the italic run\n// formatting, not just the word \"Genesis\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Find candidate paragraphs with the exact text \"Genesis\" and load their\n// font.italic so we can disambiguate the target paragraph.\nconst candidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"Genesis\") {\n    p.font.load(\"italic\");\n    candidates.push(p);\n  }\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of candidates) {\n  if (p.font.italic === true && p.style !== \"Heading 2\") {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the paragraph that holds the italic \"Genesis\" run directly\n# following the \"GEN\" (Heading 2) paragraph. There is another, unrelated\n# \"Genesis\" paragraph further down (a Heading 2 section title) that must be\n# left untouched, so we disambiguate on exact text + italic formatting +\n# style (not just the word \"Genesis\").\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $txt = $r.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"Genesis\" -and $r.Font.Italic -eq -1 -and $p.Style.NameLocal -ne \"Heading 2\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
